$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: replace "Hell 90 / gesättigt, mittel-hell" block with the new
#     "Hell 100, Sätt 100 / gesättigt, hell" fully-saturated colour block.
#     D4 (#ffa500, orange) is unchanged.
$ws.Range("C4").Value = "#ff0000"
$ws.Range("E4").Value = "#ffff00"
$ws.Range("F4").Value = "#00ff00"
$ws.Range("G4").Value = "#00ffff"
$ws.Range("H4").Value = "#0000ff"
$ws.Range("I4").Value = "#b600ff"
$ws.Range("A4").Value = "Hell 100, Sätt 100"
$ws.Range("B4").Value = "gesättigt, hell"

# --- Row 5 ("Hell 50, Sätt 50" / "mittlere Sättigung, dunkel") is unchanged.

# --- Row 6 ("Hell 100, Sätt 50" / "mittlere Sättigung, hell") moves down to
#     row 8, leaving rows 6-7 empty.
$ws.Range("A8").Value = $ws.Range("A6").Value()
$ws.Range("B8").Value = $ws.Range("B6").Value()
$ws.Range("C8").Value = $ws.Range("C6").Value()
$ws.Range("D8").Value = $ws.Range("D6").Value()
$ws.Range("E8").Value = $ws.Range("E6").Value()
$ws.Range("F8").Value = $ws.Range("F6").Value()
$ws.Range("G8").Value = $ws.Range("G6").Value()
$ws.Range("H8").Value = $ws.Range("H6").Value()
$ws.Range("I8").Value = $ws.Range("I6").Value()

$ws.Range("A6:I6").Clear()

# --- Update the active selection to match the post-edit workbook state.
$ws.Range("C11").Select()
